# Update rfuse ext4 summary results with new benchmark numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 - randread_128k
$ws.Range("B3").Value = 1052
$ws.Range("C3").Value = 2503
$ws.Range("D3").Value = 4873
$ws.Range("E3").Value = 8953
$ws.Range("F3").Value = 10800
$ws.Range("G3").Value = 13000

# Row 8 - randread_4k
$ws.Range("B8").Value = 10500
$ws.Range("C8").Value = 14400
$ws.Range("D8").Value = 28700
$ws.Range("E8").Value = 57500
$ws.Range("F8").Value = 107000
$ws.Range("G8").Value = 175000

# Row 13 - randwrite_128k
$ws.Range("B13").Value = 11800
$ws.Range("C13").Value = 15800
$ws.Range("D13").Value = 18700
$ws.Range("E13").Value = 20600
$ws.Range("F13").Value = 21600
$ws.Range("G13").Value = 21400

# Row 18 - randwrite_4k
$ws.Range("B18").Value = 366000
$ws.Range("C18").Value = 431000
$ws.Range("D18").Value = 542000
$ws.Range("E18").Value = 630000
$ws.Range("F18").Value = 612000
$ws.Range("G18").Value = 633000

# Row 23 - read_128k
$ws.Range("B23").Value = 3282
$ws.Range("C23").Value = 5876
$ws.Range("D23").Value = 8799
$ws.Range("E23").Value = 8489
$ws.Range("F23").Value = 7774
$ws.Range("G23").Value = 7627

# Row 28 - read_4k
$ws.Range("B28").Value = 102000
$ws.Range("C28").Value = 163000
$ws.Range("D28").Value = 265000
$ws.Range("E28").Value = 279000
$ws.Range("F28").Value = 274000
$ws.Range("G28").Value = 461000

# Row 33 - write_128k (only C, D, E change)
$ws.Range("C33").Value = 12100
$ws.Range("D33").Value = 13500
$ws.Range("E33").Value = 13800

# Row 38 - write_4k
$ws.Range("B38").Value = 272000
$ws.Range("C38").Value = 350000
$ws.Range("D38").Value = 415000
$ws.Range("E38").Value = 425000
$ws.Range("F38").Value = 446000
$ws.Range("G38").Value = 468000
